$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1426.6666
$ws.Range("I40").Value = 1406.8572
$ws.Range("K40").Value = 1406.8572
$ws.Range("M40").Value = -1231.8572

$ws.Range("H53").Value = 707.5
$ws.Range("I53").Value = 598.7143
$ws.Range("J53").Value = 859.8
$ws.Range("K53").Value = 598.7143
$ws.Range("L53").Value = 859.8
$ws.Range("M53").Value = 38.28570000000002
$ws.Range("N53").Value = -2133.8

$ws.Range("H80").Value = 883.6667
$ws.Range("J80").Value = 525
$ws.Range("L80").Value = 1575
$ws.Range("N80").Value = -3571

$ws.Range("H83").Value = 883.6667
$ws.Range("J83").Value = 525
$ws.Range("L83").Value = 4725
$ws.Range("N83").Value = -14709

$ws.Range("H86").Value = 1474.75
$ws.Range("I86").Value = 1950
$ws.Range("J86").Value = 999.5
$ws.Range("K86").Value = 1950
$ws.Range("L86").Value = 999.5
$ws.Range("M86").Value = -827
$ws.Range("N86").Value = -3245.5

$ws.Range("H88").Value = 1924.8
$ws.Range("I88").Value = 1933.3334
$ws.Range("J88").Value = 1912
$ws.Range("K88").Value = 1933.3334
$ws.Range("L88").Value = 1912
$ws.Range("M88").Value = -1527.3334
$ws.Range("N88").Value = -2724

$ws.Range("H89").Value = 1474.75
$ws.Range("I89").Value = 1950
$ws.Range("J89").Value = 999.5
$ws.Range("K89").Value = 9750
$ws.Range("L89").Value = 4997.5
$ws.Range("M89").Value = -4134
$ws.Range("N89").Value = -16229.5

$ws.Range("H91").Value = 1924.8
$ws.Range("I91").Value = 1933.3334
$ws.Range("J91").Value = 1912
$ws.Range("K91").Value = 1933.3334
$ws.Range("L91").Value = 1912
$ws.Range("M91").Value = -529.3334
$ws.Range("N91").Value = -4720

$ws.Range("H103").Value = 1130.3334
$ws.Range("J103").Value = 1130.3334
$ws.Range("L103").Value = 3391.0002
$ws.Range("N103").Value = -4563.0002

$ws.Range("H129").Value = 1572.3334
$ws.Range("J129").Value = 1844.6666
$ws.Range("L129").Value = 5533.9998
$ws.Range("N129").Value = -15533.9998

$ws.Range("H137").Value = 2225
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 464
$ws.Range("J16").Value = 600
$ws.Range("L16").Value = 600
$ws.Range("N16").Value = -1174

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 17399.8
$ws.Range("J95").Value = 17399.8
$ws.Range("L95").Value = 17399.8
$ws.Range("N95").Value = -22891.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -650
$ws.Range("N22").ClearContents()

$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20368

$ws.Range("H88").Value = 7332.5
$ws.Range("J88").Value = 7332.5
$ws.Range("L88").Value = 7332.5
$ws.Range("N88").Value = -8144.5

$ws.Range("H91").Value = 7332.5
$ws.Range("J91").Value = 7332.5
$ws.Range("L91").Value = 7332.5
$ws.Range("N91").Value = -10140.5

$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 20000
$ws.Range("L101").Value = 20000
$ws.Range("N101").Value = -26490

$ws.Range("H107").Value = 2332
$ws.Range("I107").Value = 2499.5
$ws.Range("J107").Value = 1997
$ws.Range("K107").Value = 2499.5
$ws.Range("L107").Value = 1997
$ws.Range("M107").Value = -579.5
$ws.Range("N107").Value = -5837

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 37.333332
$ws.Range("J2").Value = 32.84375
$ws.Range("L2").Value = 197.0625
$ws.Range("N2").Value = -423.0625

$ws.Range("H19").Value = 530.5
$ws.Range("I19").Value = 22
$ws.Range("J19").Value = 700
$ws.Range("K19").Value = 66
$ws.Range("L19").Value = 2100
$ws.Range("M19").Value = 108
$ws.Range("N19").Value = -2448

$ws.Range("H131").Value = 1950
$ws.Range("J131").Value = 1950
$ws.Range("L131").Value = 5850
$ws.Range("N131").Value = -15930

$ws.Range("H139").Value = 1584.75
$ws.Range("I139").Value = 1584.75
$ws.Range("K139").Value = 4754.25
$ws.Range("M139").Value = 385.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7340.909
$ws.Range("I80").Value = 2342
$ws.Range("K80").Value = 2342
$ws.Range("M80").Value = -1344

$ws.Range("H83").Value = 7340.909
$ws.Range("I83").Value = 2342
$ws.Range("K83").Value = 11710
$ws.Range("M83").Value = -6718

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2316.7856
$ws.Range("I7").Value = 2033.4615
$ws.Range("K7").Value = 2033.4615
$ws.Range("M7").Value = -1921.4615

$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 666.6667
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("N46").Value = -1376

$ws.Range("H99").Value = 40000
$ws.Range("I99").Value = 40000
$ws.Range("K99").Value = 40000
$ws.Range("M99").Value = -37005

$ws.Range("H126").Value = 2316.7856
$ws.Range("I126").Value = 2033.4615
$ws.Range("K126").Value = 6100.3845
$ws.Range("M126").Value = -3630.3845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 65000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 65000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H81").Value = 914.8333
$ws.Range("I81").Value = 897.8
$ws.Range("K81").Value = 1795.6
$ws.Range("M81").Value = -734.5999999999999

$ws.Range("H84").Value = 914.8333
$ws.Range("I84").Value = 897.8
$ws.Range("K84").Value = 8978
$ws.Range("M84").Value = -3674

$ws.Range("H100").Value = 374.75
$ws.Range("I100").Value = 374.75
$ws.Range("K100").Value = 749.5
$ws.Range("M100").Value = -208.5
